$d = $word.ActiveDocument

$replacements = @(
    @("64÷9=7, 1", "42÷4=10, 2"),
    @("98÷6=16, 2", "56÷5=11, 1"),
    @("11÷6=1, 5", "14÷3=4, 2"),
    @("39÷7=5, 4", "45÷7=6, 3"),
    @("24÷2=12, 0", "88÷5=17, 3"),
    @("52÷8=6, 4", "88÷2=44, 0"),
    @("50÷9=5, 5", "48÷2=24, 0"),
    @("64÷3=21, 1", "38÷5=7, 3"),
    @("29÷3=9, 2", "62÷7=8, 6"),
    @("88÷8=11, 0", "81÷8=10, 1"),
    @("48÷6=8, 0", "10÷9=1, 1"),
    @("70÷8=8, 6", "82÷4=20, 2"),
    @("71÷2=35, 1", "77÷9=8, 5"),
    @("53÷8=6, 5", "78÷7=11, 1"),
    @("90÷6=15, 0", "19÷9=2, 1"),
    @("87÷6=14, 3", "22÷2=11, 0"),
    @("49÷5=9, 4", "39÷6=6, 3"),
    @("67÷4=16, 3", "35÷4=8, 3"),
    @("74÷8=9, 2", "66÷5=13, 1"),
    @("56÷6=9, 2", "55÷3=18, 1"),
    @("18÷3=6, 0", "17÷8=2, 1"),
    @("86÷9=9, 5", "29÷6=4, 5"),
    @("99÷8=12, 3", "72÷8=9, 0"),
    @("58÷9=6, 4", "21÷8=2, 5"),
    @("92÷9=10, 2", "50÷9=5, 5")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
